$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.310.47'
$ws.Range("E2").Value = '  -0.49%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.843.02'
$ws.Range("E3").Value = '  -0.49%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.9985'
$ws.Range("E4").Value = '  -0.35%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '240.59'
$ws.Range("E5").Value = '  +0.12%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '0.6285'
$ws.Range("E6").Value = '  -0.15%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.9991'
$ws.Range("E7").Value = '  -0.29%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.07450'
$ws.Range("E8").Value = '  -2.67%  '

$ws.Range("E9").Value = '  -0.71%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '24.37'
$ws.Range("E10").Value = '  -1.87%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07716'
$ws.Range("E11").Value = '  -0.41%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '1.842.18'
$ws.Range("E12").Value = '  -2.44%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '4.992'
$ws.Range("E13").Value = '  -0.87%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.6781'
$ws.Range("E14").Value = '  -0.49%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.00001015'
$ws.Range("E15").Value = '  -4.55%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '82.01'
$ws.Range("E16").Value = '  -1.77%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '6.153'
$ws.Range("E17").Value = '  -0.70%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '29.318.84'
$ws.Range("E18").Value = '  -0.74%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '228.37'
$ws.Range("E19").Value = '  -0.27%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '12.29'
$ws.Range("E20").Value = '  -0.44%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.9994'
$ws.Range("E21").Value = '  -0.26%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.383'
$ws.Range("E22").Value = '  -1.09%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.9987'
$ws.Range("E23").Value = '  -0.28%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '158.61'
$ws.Range("E24").Value = '  +0.64%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.1374'
$ws.Range("E25").Value = '  -0.85%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.411'
$ws.Range("E26").Value = '  -0.25%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '17.54'
$ws.Range("E27").Value = '  -1.14%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.06493'
$ws.Range("E28").Value = '  +15.74%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.395'
$ws.Range("E29").Value = '  +0.23%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.475'
$ws.Range("E30").Value = '  +0.70%  '

$ws.Range("E31").Value = '  -1.38%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.047'
$ws.Range("E32").Value = '  -0.46%  '

$ws.Range("E33").Value = '  -1.52%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.141'
$ws.Range("E34").Value = '  -2.12%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.6982'
$ws.Range("E35").Value = '  -0.29%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.580'
$ws.Range("E36").Value = '  -0.47%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.835'
$ws.Range("E37").Value = '  +3.29%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.250.37'
$ws.Range("E38").Value = '  +1.53%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.01813'
$ws.Range("E39").Value = '  +0.43%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.513'
$ws.Range("E40").Value = '  +0.62%  '

$ws.Range("E41").Value = '  +0.14%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9986'
$ws.Range("E42").Value = '  -0.34%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '2.003.79'
$ws.Range("E43").Value = '  -16.20%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '101.12'
$ws.Range("E44").Value = '  -0.86%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '66.26'
$ws.Range("E45").Value = '  +0.24%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '7.044'
$ws.Range("E46").Value = '  -2.22%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.1168'
$ws.Range("E47").Value = '  +1.18%  '

$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.035'
$ws.Range("E48").Value = '  -0.07%  '

$ws.Range("B49").Value = 'TheSandbox'
$ws.Range("C49").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.3943'
$ws.Range("E49").Value = '  -2.25%  '

$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.670'
$ws.Range("E50").Value = '  -0.73%  '

$ws.Range("B51").Value = 'BabyDogeCoin'
$ws.Range("C51").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.00000000114'
$ws.Range("E51").Value = '  -3.51%  '
